# Commit: "Removed level specification in instructions; fixed small errors"
#
# On the "Trend_instructions" sheet, the four columns that spelled out the
# fuel/sector aggregation levels (L1_agg_fuel, L2_CEDS_fuel, L3_agg_sector,
# L4_CEDS_sector -- columns F:I) are removed. The following columns
# (override_normalization, use_as_trend, match_year) shift left to take
# their place (F:H). Deleting the entire columns (rather than clearing
# cell contents) naturally shifts the remaining columns, updates the
# sheet dimension, and drops the now-unused shared strings for the
# removed level-spec labels.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Trend_instructions")

$ws.Range("F1:I1").EntireColumn.Delete()

# The user's selection ended up on J15 after the edit.
[void]$ws.Range("J15").Select()
